$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---

# Cells whose new price text is ambiguous with a plain number (e.g. "217.00")
# need the column pre-set to Text format so Excel keeps the literal digits
# (incl. trailing/leading zeros) instead of silently converting to a Number.
$textForceCells = @('D5', 'D6', 'D7', 'D8', 'D9', 'D10', 'D11', 'D13', 'D15', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D32', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D41', 'D44', 'D47', 'D48', 'D49', 'D51')
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '217.00'
$ws.Range('D6').Value = '0.5137'
$ws.Range('D7').Value = '1.004'
$ws.Range('D8').Value = '0.2639'
$ws.Range('D9').Value = '0.06266'
$ws.Range('D10').Value = '20.74'
$ws.Range('D11').Value = '0.07749'
$ws.Range('D13').Value = '4.448'
$ws.Range('D15').Value = '0.5436'
$ws.Range('D20').Value = '4.629'
$ws.Range('D21').Value = '192.52'
$ws.Range('D22').Value = '10.07'
$ws.Range('D23').Value = '6.014'
$ws.Range('D25').Value = '140.01'
$ws.Range('D26').Value = '0.1222'
$ws.Range('D27').Value = '7.227'
$ws.Range('D29').Value = '1.434'
$ws.Range('D30').Value = '0.05946'
$ws.Range('D32').Value = '3.564'
$ws.Range('D34').Value = '1.600'
$ws.Range('D35').Value = '0.9653'
$ws.Range('D36').Value = '2.423'
$ws.Range('D37').Value = '2.773'
$ws.Range('D38').Value = '0.5649'
$ws.Range('D39').Value = '0.01590'
$ws.Range('D40').Value = '5.961'
$ws.Range('D41').Value = '0.8564'
$ws.Range('D44').Value = '100.22'
$ws.Range('D47').Value = '56.63'
$ws.Range('D48').Value = '1.007'
$ws.Range('D49').Value = '7.999'
$ws.Range('D51').Value = '1.455'

# Cells whose new price text already cannot be parsed as a plain number
# (multiple "." separators, or subscript-digit notation) -> safe to set directly.
$ws.Range('D2').Value = '26.175.85'
$ws.Range('D3').Value = '1.660.01'
$ws.Range('D12').Value = '1.666.13'
$ws.Range('D14').Value = '1.887.15'
$ws.Range('D16').Value = '0.0₅8099'
$ws.Range('D18').Value = '26.196.94'
$ws.Range('D43').Value = '1.012.24'
$ws.Range('D45').Value = '1.801.68'
$ws.Range('D46').Value = '0.0₈111'

# --- Volume(1h) column (E) updates ---

$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('E3').Value = '  -0.86%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('E5').Value = '  -1.57%  '
$ws.Range('E6').Value = '  -3.27%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E8').Value = '  -1.86%  '
$ws.Range('E9').Value = '  -1.92%  '
$ws.Range('E10').Value = '  -4.53%  '
$ws.Range('E11').Value = '  -0.61%  '
$ws.Range('E12').Value = '  -0.88%  '
$ws.Range('E13').Value = '  -0.98%  '
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('E17').Value = '  -1.12%  '
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  -2.83%  '
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('E23').Value = '  -4.78%  '
$ws.Range('E25').Value = '  +1.21%  '
$ws.Range('E26').Value = '  -4.18%  '
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('E29').Value = '  +0.30%  '
$ws.Range('E30').Value = '  -5.26%  '
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('E33').Value = '  -4.77%  '
$ws.Range('E34').Value = '  -5.38%  '
$ws.Range('E35').Value = '  -4.39%  '
$ws.Range('E36').Value = '  +0.15%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  -8.22%  '
$ws.Range('E39').Value = '  -1.81%  '
$ws.Range('E40').Value = '  -2.79%  '
$ws.Range('E41').Value = '  -0.57%  '
$ws.Range('E42').Value = '  +0.26%  '
$ws.Range('E43').Value = '  -7.30%  '
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('E47').Value = '  -3.67%  '
$ws.Range('E48').Value = '  +0.28%  '
$ws.Range('E49').Value = '  -2.35%  '
$ws.Range('E51').Value = '  -4.26%  '

